# Update (Removed Auto Arima)
# Updates the "Forecast Comparison" sheet's forecast columns (D: Amazon Mean
# Forecast, E: Amazon P70 Forecast, F: Amazon P80 Forecast, G: Amazon P90
# Forecast, and the recalculated C: Prophet Forecast) for weeks 1-16 (rows
# 2-17) now that the Auto-ARIMA model has been removed from the forecast
# blend, and refreshes the dependent totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$fc = $wb.Worksheets.Item("Forecast Comparison")

# row -> (C: Prophet, D: Amazon Mean, E: Amazon P70, F: Amazon P80, G: Amazon P90)
$rows = @{
    2  = @(45, 51, 60, 68, 81)
    3  = @(44, 41, 49, 57, 69)
    4  = @(44, 36, 43, 50, 61)
    5  = @(44, 36, 44, 51, 62)
    6  = @(41, 36, 44, 51, 62)
    7  = @(38, 37, 44, 51, 62)
    8  = @(37, 37, 45, 52, 64)
    9  = @(39, 37, 45, 52, 64)
    10 = @(43, 36, 43, 51, 62)
    11 = @(45, 36, 43, 51, 63)
    12 = @(42, 36, 44, 52, 64)
    13 = @(40, 37, 45, 54, 68)
    14 = @(37, 36, 43, 51, 63)
    15 = @(38, 35, 43, 52, 65)
    16 = @(37, 35, 42, 50, 62)
    17 = @(33, 34, 41, 49, 61)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $fc.Range("C$r").Value = $vals[0]
    $fc.Range("D$r").Value = $vals[1]
    $fc.Range("E$r").Value = $vals[2]
    $fc.Range("F$r").Value = $vals[3]
    $fc.Range("G$r").Value = $vals[4]
}

$sum = $wb.Worksheets.Item("Summary")

# These cells hold numeric-looking text (matching the rest of column B on
# this sheet), so force text entry with a leading apostrophe to avoid
# Excel auto-converting them to the Number type.
$sum.Range("B9").Value = "'647"
$sum.Range("B10").Value = "'332"
$sum.Range("B11").Value = "'177"
$sum.Range("B14").Value = "'33"

Write-Output "Updated forecast columns and summary totals after removing Auto-ARIMA"
